$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated imputed values (RandomForest algorithm rerun) for the
# affected cells, as reflected by the upstream result data refresh.
$ws.Range("D7").Value = -7.409000000000003
$ws.Range("B8").Value = 5.750799999999998
$ws.Range("B10").Value = 5.156700000000001
$ws.Range("B12").Value = 4.944900000000001
$ws.Range("D15").Value = -8.015599999999997
$ws.Range("B18").Value = 5.532299999999997
$ws.Range("D18").Value = -8.292099999999991
$ws.Range("D20").Value = -7.742899999999995
$ws.Range("D29").Value = -7.474999999999999
$ws.Range("D30").Value = -7.383099999999998
$ws.Range("D31").Value = -8.548699999999998
$ws.Range("B37").Value = 8.881
$ws.Range("D40").Value = -7.982499999999995
$ws.Range("D50").Value = -8.199699999999998
$ws.Range("B55").Value = 5.110999999999997
$ws.Range("B68").Value = 4.649299999999998
$ws.Range("D68").Value = -7.112799999999994
$ws.Range("D76").Value = -7.390699999999998
$ws.Range("B77").Value = 9.133100000000004
$ws.Range("B78").Value = 9.691900000000002
$ws.Range("B81").Value = 5.111200000000003
$ws.Range("B82").Value = 5.687800000000003
$ws.Range("D87").Value = -7.927999999999995
$ws.Range("D88").Value = -7.544399999999996
$ws.Range("D96").Value = -7.375100000000002
$ws.Range("D98").Value = -8.271600000000001
$ws.Range("D101").Value = -7.842199999999998
$ws.Range("D102").Value = -7.783199999999999
